$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F header value (new shared string "LEGAL.LETTER:1")
$ws.Range("F1").Value = "LEGAL.LETTER:1"

# Reflect the saved selection state after the edit
$ws.Range("H20").Select()
